$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 194, shifting the existing "row-194-onward" block down.
$ws.Rows.Item(194).Insert()
$ws.Rows.Item(194).Insert()

# Populate the two newly inserted rows with the new weekly observation (2021-09-10).
$ws.Range('A194').Value = 8
$ws.Range('B194').Value = 'Terminal La Palmera de La Serena'
$ws.Range('C194').Value = 'Coquimbo'
$ws.Range('D194').Value = 44449
$ws.Range('E194').Value = 4
$ws.Range('F194').Value = 100112017
$ws.Range('G194').Value = 'Apio'
$ws.Range('H194').Value = 'Americana (o)'
$ws.Range('I194').Value = 'Primera'
$ws.Range('J194').Value = 2000
$ws.Range('K194').Value = 7500
$ws.Range('L194').Value = 8000
$ws.Range('M194').Value = 7750
$ws.Range('N194').Value = '$/docena de matas'
$ws.Range('O194').Value = 'Provincia del Elquí'
$ws.Range('P194').Value = 1292
$ws.Range('Q194').Value = 6
$ws.Range('R194').Value = 'Hortaliza'
$ws.Range('A195').Value = 8
$ws.Range('B195').Value = 'Terminal La Palmera de La Serena'
$ws.Range('C195').Value = 'Coquimbo'
$ws.Range('D195').Value = 44449
$ws.Range('E195').Value = 4
$ws.Range('F195').Value = 100112017
$ws.Range('G195').Value = 'Apio'
$ws.Range('H195').Value = 'Americana (o)'
$ws.Range('I195').Value = 'Segunda'
$ws.Range('J195').Value = 1400
$ws.Range('K195').Value = 6500
$ws.Range('L195').Value = 7000
$ws.Range('M195').Value = 6750
$ws.Range('N195').Value = '$/docena de matas'
$ws.Range('O195').Value = 'Provincia del Elquí'
$ws.Range('P195').Value = 1125
$ws.Range('Q195').Value = 6
$ws.Range('R195').Value = 'Hortaliza'
